# Fix double dashes / stray leading spaces caused by PowerPoint autoformat
# in job submission script examples (scheduling_jobs deck).
#
# Edits are applied per-run via TextRange.Characters(start, length) so that
# the existing run (formatting) boundaries are preserved exactly; within a
# paragraph, edits are always issued from the highest start offset to the
# lowest so that previously-computed (pre-edit) offsets stay valid.

function Set-RunText {
    param($ParaRange, [int]$Start, [int]$Length, [string]$NewText)
    $rng = $ParaRange.Characters($Start, $Length)
    $rng.Text = $NewText
}

$p = $ppt.ActivePresentation

# --- Slide 18: "Example job script" --------------------------------------
$s18 = $p.Slides.Item(18)
$tr18 = $s18.Shapes.Item(3).TextFrame.TextRange

Set-RunText $tr18.Paragraphs(1, 1) 1 12 "#!/bin/bash"
Set-RunText $tr18.Paragraphs(3, 1) 1 14 "## Directives"
Set-RunText $tr18.Paragraphs(4, 1) 1 11 "#SBATCH --"
Set-RunText $tr18.Paragraphs(5, 1) 1 52 "#SBATCH --time=0:01:00              `t# Max run time"

$para6_18 = $tr18.Paragraphs(6, 1)
Set-RunText $para6_18 28 32 "      `t       # Specify Alpine CPU node"
Set-RunText $para6_18 1 21 "#SBATCH --partition="

$para7_18 = $tr18.Paragraphs(7, 1)
Set-RunText $para7_18 15 38 "=normal                   # Specify QoS"
Set-RunText $para7_18 1 11 "#SBATCH --"

Set-RunText $tr18.Paragraphs(8, 1) 1 24 "#SBATCH --output=test_%"
Set-RunText $tr18.Paragraphs(10, 1) 1 12 "## Software"
Set-RunText $tr18.Paragraphs(11, 1) 1 67 "module purge                          # Purge all existing modules"
Set-RunText $tr18.Paragraphs(13, 1) 1 17 "## User commands"
Set-RunText $tr18.Paragraphs(14, 1) 1 37 'echo "This is a test submitted by $USER" '

# --- Slide 24: "Software job script example" ------------------------------
$s24 = $p.Slides.Item(24)
$tr24 = $s24.Shapes.Item(3).TextFrame.TextRange

Set-RunText $tr24.Paragraphs(1, 1) 1 12 "#!/bin/bash"
Set-RunText $tr24.Paragraphs(3, 1) 1 14 "## Directives"
Set-RunText $tr24.Paragraphs(4, 1) 1 11 "#SBATCH --"
Set-RunText $tr24.Paragraphs(5, 1) 1 52 "#SBATCH --time=0:01:00              `t# Max run time"

$para6_24 = $tr24.Paragraphs(6, 1)
Set-RunText $para6_24 28 32 "      `t       # Specify Alpine CPU node"
Set-RunText $para6_24 1 21 "#SBATCH --partition="

$para7_24 = $tr24.Paragraphs(7, 1)
Set-RunText $para7_24 15 38 "=normal                   # Specify QoS"
Set-RunText $para7_24 1 11 "#SBATCH --"

Set-RunText $tr24.Paragraphs(8, 1) 1 24 "#SBATCH --output=test_%"
Set-RunText $tr24.Paragraphs(10, 1) 1 12 "## Software"
Set-RunText $tr24.Paragraphs(11, 1) 1 69 "module purge                            # Purge all existing modules"
Set-RunText $tr24.Paragraphs(12, 1) 1 56 "module load anaconda                    # Load Anaconda"
Set-RunText $tr24.Paragraphs(13, 1) 1 1 ""
Set-RunText $tr24.Paragraphs(15, 1) 1 21 "## Run Python script"
Set-RunText $tr24.Paragraphs(16, 1) 1 8 "python "

# --- Slide 26: "Running an interactive job" --------------------------------
$s26 = $p.Slides.Item(26)
$tr26 = $s26.Shapes.Item(4).TextFrame.TextRange

$para1_26 = $tr26.Paragraphs(1, 1)
Set-RunText $para1_26 41 24 "=compile --time=00:10:00"
Set-RunText $para1_26 36 2 " --"
